$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '25.854.82'
Set-TextValue $ws.Range('E2') '  -1.09%  '
Set-TextValue $ws.Range('D3') '1.634.83'
Set-TextValue $ws.Range('E3') '  -0.97%  '
Set-TextValue $ws.Range('D4') '1.001'
Set-TextValue $ws.Range('E4') '  -0.31%  '
Set-TextValue $ws.Range('D5') '214.69'
Set-TextValue $ws.Range('E5') '  -0.24%  '
Set-TextValue $ws.Range('D6') '0.5015'
Set-TextValue $ws.Range('E6') '  -1.90%  '
Set-TextValue $ws.Range('D7') '1.001'
Set-TextValue $ws.Range('E7') '  -0.35%  '
Set-TextValue $ws.Range('D8') '0.2559'
Set-TextValue $ws.Range('E8') '  -1.02%  '
Set-TextValue $ws.Range('D9') '0.06371'
Set-TextValue $ws.Range('E9') '  -0.91%  '
Set-TextValue $ws.Range('D10') '19.67'
Set-TextValue $ws.Range('E10') '  -1.30%  '
Set-TextValue $ws.Range('D11') '0.07711'
Set-TextValue $ws.Range('E11') '  -0.99%  '
Set-TextValue $ws.Range('D12') '1.662.21'
Set-TextValue $ws.Range('E12') '  +0.64%  '
Set-TextValue $ws.Range('D13') '4.260'
Set-TextValue $ws.Range('E13') '  -0.50%  '
Set-TextValue $ws.Range('D14') '1.862.19'
Set-TextValue $ws.Range('E14') '  -0.86%  '
Set-TextValue $ws.Range('D15') '0.5443'
Set-TextValue $ws.Range('E15') '  -1.24%  '
Set-TextValue $ws.Range('D16') '0.0₅7891'
Set-TextValue $ws.Range('E16') '  -1.38%  '
Set-TextValue $ws.Range('D17') '64.16'
Set-TextValue $ws.Range('E17') '  +0.17%  '
Set-TextValue $ws.Range('D18') '25.869.90'
Set-TextValue $ws.Range('E18') '  -1.04%  '
Set-TextValue $ws.Range('D19') '1.002'
Set-TextValue $ws.Range('E19') '  -0.27%  '
Set-TextValue $ws.Range('D20') '202.84'
Set-TextValue $ws.Range('E20') '  -3.85%  '
Set-TextValue $ws.Range('D21') '4.365'
Set-TextValue $ws.Range('E21') '  -0.59%  '
Set-TextValue $ws.Range('D22') '9.883'
Set-TextValue $ws.Range('E22') '  -1.75%  '
Set-TextValue $ws.Range('D23') '5.969'
Set-TextValue $ws.Range('E23') '  -1.18%  '
Set-TextValue $ws.Range('D24') '1.002'
Set-TextValue $ws.Range('E24') '  -0.27%  '
Set-TextValue $ws.Range('D25') '1.925'
Set-TextValue $ws.Range('E25') '  +10.00%  '
Set-TextValue $ws.Range('D26') '141.10'
Set-TextValue $ws.Range('E26') '  -1.89%  '
Set-TextValue $ws.Range('D27') '0.1134'
Set-TextValue $ws.Range('E27') '  -3.39%  '
Set-TextValue $ws.Range('D28') '15.69'
Set-TextValue $ws.Range('E28') '  -0.70%  '
Set-TextValue $ws.Range('D29') '6.707'
Set-TextValue $ws.Range('E29') '  -3.86%  '
Set-TextValue $ws.Range('D30') '1.240'
Set-TextValue $ws.Range('E30') '  -0.08%  '
Set-TextValue $ws.Range('D31') '0.04928'
Set-TextValue $ws.Range('E31') '  -4.02%  '
Set-TextValue $ws.Range('D32') '3.269'
Set-TextValue $ws.Range('E32') '  -2.19%  '
Set-TextValue $ws.Range('D33') '3.182'
Set-TextValue $ws.Range('E33') '  -1.02%  '
Set-TextValue $ws.Range('D34') '1.537'
Set-TextValue $ws.Range('E34') '  -1.29%  '
Set-TextValue $ws.Range('D35') '2.366'
Set-TextValue $ws.Range('E35') '  +0.73%  '
Set-TextValue $ws.Range('D36') '2.624'
Set-TextValue $ws.Range('E36') '  -4.19%  '
Set-TextValue $ws.Range('D37') '0.8914'
Set-TextValue $ws.Range('E37') '  -3.47%  '
Set-TextValue $ws.Range('D38') '1.155.92'
Set-TextValue $ws.Range('E38') '  -1.02%  '
Set-TextValue $ws.Range('D39') '0.5591'
Set-TextValue $ws.Range('E39') '  -1.72%  '
Set-TextValue $ws.Range('D40') '0.01560'
Set-TextValue $ws.Range('E40') '  -1.42%  '
Set-TextValue $ws.Range('D41') '0.9999'
Set-TextValue $ws.Range('E41') '  -0.40%  '
Set-TextValue $ws.Range('D42') '5.704'
Set-TextValue $ws.Range('E42') '  +0.85%  '
Set-TextValue $ws.Range('D43') '0.8069'
Set-TextValue $ws.Range('E43') '  -2.08%  '
Set-TextValue $ws.Range('E44') '  -0.48%  '
Set-TextValue $ws.Range('D45') '1.773.92'
Set-TextValue $ws.Range('E45') '  -0.85%  '
Set-TextValue $ws.Range('E46') '  -0.53%  '
Set-TextValue $ws.Range('D47') '0.4513'
Set-TextValue $ws.Range('E47') '  -0.80%  '
Set-TextValue $ws.Range('D48') '1.005'
Set-TextValue $ws.Range('E48') '  -0.05%  '
Set-TextValue $ws.Range('D49') '54.88'
Set-TextValue $ws.Range('E49') '  -0.96%  '
Set-TextValue $ws.Range('D50') '0.05052'
Set-TextValue $ws.Range('E50') '  -0.31%  '
Set-TextValue $ws.Range('E51') '  -0.20%  '

Write-Output "done"